# Cotações atualizadas - 2025-11-01
# Append a new row (58) with the quotes for 2025-11-01 (Excel serial date 45962),
# mirroring the formatting of the preceding data row (57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 58
$prevRow = 57

# Date column (A): keep the same date/number formatting as the row above.
$ws.Cells.Item($newRow, 1).Value = 45962
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

# Quote columns (B:E): stored as text values using comma decimal separators.
$ws.Cells.Item($newRow, 2).Value = "22,0341"
$ws.Cells.Item($newRow, 3).Value = "16,1343"
$ws.Cells.Item($newRow, 4).Value = "15,5326"
$ws.Cells.Item($newRow, 5).Value = "15,5326"
